$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression values
$ws.Range("B2").Value = 5877834119531547
$ws.Range("C2").Value = 5877834119531547
$ws.Range("D2").Value = 5877834119531547

# Row 3 - RandomForestRegressor values
$ws.Range("B3").Value = 177469956896346.4
$ws.Range("C3").Value = 159962917570449.5
$ws.Range("D3").Value = 737471502163255.2

# Row 4 - rename model and update values
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 8932460205461.195
$ws.Range("C4").Value = 9438776346173.814
$ws.Range("D4").Value = 192756897844512.3

# Row 5 - rename model and update values
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 393941520254848.1
$ws.Range("C5").Value = 487676794515036.2
$ws.Range("D5").Value = 2808045136390766
